# Update countries & provincias Spain
# Applies the 16-Jun-2020 00:28 data refresh to the "Pais" sheet:
#   - bump the "Datos actualizados..." timestamp in A1
#   - refresh case counters for several countries
#   - two countries in the (now-updated) ranking swapped places with their
#     neighbour, so the two rows involved need both their name and their
#     numbers updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp ---------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 16 de Junio de 2020 a las 00:28"

# --- Plain numeric refreshes (country order unchanged) --------------------

# Row 4: Estados Unidos
$ws.Range("B4").Value = 2181553
$ws.Range("C4").Value = 19325
$ws.Range("D4").Value = 879295
$ws.Range("E4").Value = 1183995
$ws.Range("G4").Value = 405
$ws.Range("H4").Value = 118263

# Row 5: Brasil
$ws.Range("B5").Value = 888271
$ws.Range("C5").Value = 20389
$ws.Range("E5").Value = 390744
$ws.Range("G5").Value = 570
$ws.Range("H5").Value = 43959

# Row 11: Peru
$ws.Range("B11").Value = 232992
$ws.Range("C11").Value = 3256
$ws.Range("D11").Value = 119409
$ws.Range("E11").Value = 106723
$ws.Range("G11").Value = 172
$ws.Range("H11").Value = 6860

# Row 24: Sudafrica
$ws.Range("B24").Value = 73533
$ws.Range("C24").Value = 3495
$ws.Range("D24").Value = 39867
$ws.Range("E24").Value = 32098
$ws.Range("G24").Value = 88
$ws.Range("H24").Value = 1568

# Row 90: Bulgaria
$ws.Range("B90").Value = 3341
$ws.Range("C90").Value = 51
$ws.Range("D90").Value = 1784
$ws.Range("E90").Value = 1381
$ws.Range("G90").Value = 2
$ws.Range("H90").Value = 176

# Row 161: Surinam
$ws.Range("B161").Value = 229
$ws.Range("C161").Value = 21
$ws.Range("D161").Value = 48
$ws.Range("E161").Value = 177

# Row 197: Curazao
$ws.Range("D197").Value = 19
$ws.Range("E197").Value = 2

# --- Rows 27/28: Suecia and Colombia swap rank order -----------------------
# Row 27 used to be Suecia, now becomes Colombia (with Colombia's refreshed
# numbers); row 28 used to be Colombia, now becomes Suecia (keeping Suecia's
# previous numbers, unchanged this update).
$ws.Range("A27").Value = "Colombia"
$ws.Range("B27").Value = 53063
$ws.Range("C27").Value = 2124
$ws.Range("D27").Value = 19952
$ws.Range("E27").Value = 31385
$ws.Range("G27").Value = 59
$ws.Range("H27").Value = 1726

$ws.Range("A28").Value = "Suecia"
$ws.Range("B28").Value = 52383
$ws.Range("C28").Value = 139
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 17
$ws.Range("H28").Value = 4891

# --- Rows 208-211: Santa Sede / Islas Turcas y Caicos / Seychelles /
#     Montserrat reshuffle ---------------------------------------------------
# Only "Casos activos" (D) and "Muertes" (H) differ between these four
# territories, so only the name + D/H need to move.
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1

$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1
